$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8549234867095947
$ws.Range("B1").Value = 1.238696575164795
$ws.Range("C1").Value = 2.421857833862305
$ws.Range("D1").Value = 3.770876407623291
$ws.Range("E1").Value = 1.894841432571411
